{"js": "// The document contains one table of two-digit \u00f7 one-digit division problems\n// (some rows are intentionally blank spacer rows). This edit swaps each\n// problem/answer string for a new one, in left-to-right, top-to-bottom cell\n// order. Several of the old/new strings repeat elsewhere in the table, so the\n// replacement must be positional (by cell), not a global text search/replace.\n\nconst replacements = [\n  \"18\u00f75=3, 3\", \"59\u00f72=29, 1\", \"37\u00f79=4, 1\", \"82\u00f78=10, 2\", \"41\u00f76=6, 5\",\n  \"11\u00f74=2, 3\", \"99\u00f73=33, 0\", \"85\u00f72=42, 1\", \"96\u00f77=13, 5\", \"71\u00f76=11, 5\",\n  \"27\u00f79=3, 0\", \"90\u00f79=10, 0\", \"64\u00f78=8, 0\", \"70\u00f73=23, 1\", \"84\u00f72=42, 0\",\n  \"50\u00f78=6, 2\", \"33\u00f72=16, 1\", \"59\u00f72=29, 1\", \"71\u00f76=11, 5\", \"62\u00f73=20, 2\",\n  \"80\u00f79=8, 8\", \"27\u00f78=3, 3\", \"19\u00f79=2, 1\", \"35\u00f77=5, 0\", \"10\u00f75=2, 0\"\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Load the cell collection for every row.\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\n// Load the paragraph collection for every cell (one paragraph per cell here).\nconst cellParagraphsList = [];\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    const paragraphs = cell.body.paragraphs;\n    paragraphs.load(\"items\");\n    cellParagraphsList.push(paragraphs);\n  }\n}\nawait context.sync();\n\n// Load the text of each cell's first paragraph so we can tell which cells\n// hold a division problem versus an intentionally empty spacer cell.\nconst firstParas = cellParagraphsList.map((paragraphs) => paragraphs.items[0]);\nfor (const para of firstParas) {\n  para.load(\"text\");\n}\nawait context.sync();\n\n// Walk the cells in row-major, left-to-right order (matching the order the\n// division problems appear in the source diff) and replace each non-blank\n// cell's text in turn, preserving the run's existing formatting.\nlet replacementIndex = 0;\nfor (const para of firstParas) {\n  if (para.text && para.text.trim() !== \"\") {\n    para.insertText(replacements[replacementIndex], Word.InsertLocation.replace);\n    replacementIndex++;\n  }\n}\nawait context.sync();\n", "ps1": "# The document contains one table of two-digit \u00f7 one-digit division problems.\n# This edit swaps each problem/answer string for a new one, in left-to-right,\n# top-to-bottom cell order (content rows only - some rows are intentionally\n# blank spacer rows). Because several of the old/new strings repeat elsewhere\n# in the table, the replacement must be positional (by cell), not a global\n# text search/replace.\n\n$replacements = @(\n    \"18\u00f75=3, 3\", \"59\u00f72=29, 1\", \"37\u00f79=4, 1\", \"82\u00f78=10, 2\", \"41\u00f76=6, 5\",\n    \"11\u00f74=2, 3\", \"99\u00f73=33, 0\", \"85\u00f72=42, 1\", \"96\u00f77=13, 5\", \"71\u00f76=11, 5\",\n    \"27\u00f79=3, 0\", \"90\u00f79=10, 0\", \"64\u00f78=8, 0\", \"70\u00f73=23, 1\", \"84\u00f72=42, 0\",\n    \"50\u00f78=6, 2\", \"33\u00f72=16, 1\", \"59\u00f72=29, 1\", \"71\u00f76=11, 5\", \"62\u00f73=20, 2\",\n    \"80\u00f79=8, 8\", \"27\u00f78=3, 3\", \"19\u00f79=2, 1\", \"35\u00f77=5, 0\", \"10\u00f75=2, 0\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n$index = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cellText = $cell.Range.Text -replace \"[\\r\\x07]\", \"\"\n        if ($cellText.Trim() -ne \"\") {\n            $cell.Range.Text = $replacements[$index]\n            $index++\n        }\n    }\n}\n"}
